# Generate Report for Handback
# The handback locale changes from "de-de" to "zh-tw": rename the locale
# worksheet/table, refresh the two handback timestamps on that sheet, and
# relabel the locale column on the "Overview" summary sheet.

$wb = $excel.ActiveWorkbook

# --- Locale worksheet: "de-de" -> "zh-tw" -----------------------------
$wsLocale = $wb.Worksheets.Item("de-de")
$wsLocale.Name = "zh-tw"

# Refresh the handback datetimes recorded for this run.
$wsLocale.Range("E2:E5").Value = "2016-03-11 01:03:15"
$wsLocale.Range("H2:H5").Value = "2016-03-17 04:12:10"

# Rename the locale table that lives on that worksheet to match.
$loLocale = $wsLocale.ListObjects.Item(1)
$loLocale.Name = "zh-tw"

# --- Overview worksheet: relabel the locale column ---------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# The Overview table's header row isn't shown (headerRowCount = 0), so the
# table-column name is decoupled from the cell text. Toggle the header row
# on just long enough for a cell edit to resync the column name, then
# restore the original (headerless) layout and range.
$overviewRange = $loOverview.Range
$loOverview.ShowHeaders = $true
$wsOverview.Range("B1").Value = "zh-tw"
$loOverview.ShowHeaders = $false
$loOverview.Resize($overviewRange)
